$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item("Hoja2").Name = "Questionary"
$wb.Worksheets.Item("Hoja1").Name = "Task"
